$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A18").Value = "Wil je dit even doorsturen?"
$ws.Range("B18").Value = "mailmind.test@zohomail.eu"
$ws.Range("C18").Value = "Testmail #16: Wil je dit even doorsturen?"
$ws.Range("D18").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E18").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$ws.Range("F18").Value = "2025-07-31 21:56:55"
$ws.Range("G18").Value = "Ja"
$ws.Range("H18").Value = "Ja"
$ws.Range("I18").Value = "Nee"
$ws.Range("J18").Value = "Nee"

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = 4

$ws.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D18"))
$ws.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G18"))
$ws.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H18"))
$ws.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I18"))
$ws.Range("J2:J17").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J18"))
